$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old B1 value ("tgrt") - it is replaced by new cells in column C.
$ws.Range("B1").ClearContents()

# Add the three new shared-string values in column C, rows 2-4.
$ws.Range("C2").Value = "gtyu"
$ws.Range("C3").Value = "gf"
$ws.Range("C4").Value = "fff"

# Move the active selection to B1 to match the edited view state.
[void]$ws.Range("B1").Select()
